$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: New Jersey
$ws.Range("A25").Value = "New Jersey"
$ws.Range("B25").Value = 46
$ws.Range("C25").Value = 64
$ws.Range("D25").Value = 445
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 98
$ws.Range("G25").Value = 378
$ws.Range("H25").Formula = "=445-G25"
$ws.Range("I25").Value = 49
$ws.Range("J25").Value = 434
$ws.Range("K25").Formula = "=445-J25"
$ws.Range("L25").Value = 81
$ws.Range("M25").Value = 289
$ws.Range("N25").Formula = "=445-M25"
$ws.Range("O25").Value = 152
$ws.Range("P25").Value = 222
$ws.Range("Q25").Formula = "=445-P25"

# Row 26: New York
$ws.Range("A26").Value = "New York"
$ws.Range("B26").Value = 45
$ws.Range("C26").Value = 74
$ws.Range("D26").Value = 451
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 91
$ws.Range("G26").Value = 371
$ws.Range("H26").Formula = "=451-G26"
$ws.Range("I26").Value = 45
$ws.Range("J26").Value = 439
$ws.Range("K26").Formula = "=451-J26"
$ws.Range("L26").Value = 79
$ws.Range("M26").Value = 255
$ws.Range("N26").Formula = "=451-M26"
$ws.Range("O26").Value = 159
$ws.Range("P26").Value = 230
$ws.Range("Q26").Formula = "=451-P26"

# Row 27: North Carolina
$ws.Range("A27").Value = "North Carolina"
$ws.Range("B27").Value = 71
$ws.Range("C27").Value = 255
$ws.Range("D27").Value = 352
$ws.Range("E27").Formula = "=403-D27"
$ws.Range("F27").Value = 22
$ws.Range("G27").Value = 349
$ws.Range("H27").Formula = "=403-G27"
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 397
$ws.Range("K27").Formula = "=403-J27"
$ws.Range("L27").Value = 50
$ws.Range("M27").Value = 179
$ws.Range("N27").Formula = "=403-M27"
$ws.Range("O27").Value = 59
$ws.Range("P27").Value = 206
$ws.Range("Q27").Formula = "=403-P27"

# Row 28: Ohio
$ws.Range("A28").Value = "Ohio"
$ws.Range("B28").Value = 33
$ws.Range("C28").Value = 127
$ws.Range("D28").Value = 473
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 62
$ws.Range("G28").Value = 421
$ws.Range("H28").Formula = "=473-G28"
$ws.Range("I28").Value = 27
$ws.Range("J28").Value = 466
$ws.Range("K28").Formula = "=473-J28"
$ws.Range("L28").Value = 53
$ws.Range("M28").Value = 290
$ws.Range("N28").Formula = "=473-M28"
$ws.Range("O28").Value = 201
$ws.Range("P28").Value = 240
$ws.Range("Q28").Formula = "=473-P28"

# Extend the shared "control number" formula down through the new rows in
# one bulk assignment so it is stored as a shared formula group, matching
# how Excel's fill-down records R3:R24 -> R3:R28.
$ws.Range("R25:R28").Formula = "=SUM(C25+F25+I25+L25+O25)-AVERAGE(SUM(E25+D25),SUM(G25+H25),SUM(J25+K25), SUM(M25:N25),SUM(P25:Q25))"

$ws.Range("Q28").Select()
